{"js": "// Update the SDK compatibility line:\n//   \"This sample is compatible with the Windows 10 Creators Update SDK (15063)\"\n// becomes:\n//   \"This sample is compatible with the Windows 10 Fall Creators Update SDK (16299)\"\n// (commit message: \"Updated for Windows 10 Fall Creators Update SDK (16299)\")\n//\n// Word records this as the user selecting just the version substring and\n// retyping it, which both splits the run at that point and relocates the\n// \"_GoBack\" last-edit bookmark to right after the new text. We replicate\n// both effects here.\n\nconst doc = context.document;\n\nconst results = doc.body.search(\"Windows 10 Creators Update SDK (15063)\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text 'Windows 10 Creators Update SDK (15063)' not found.\");\n}\n\nconst target = results.items[0];\ntarget.insertText(\"Windows 10 Fall Creators Update SDK (16299)\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark (Word's marker for the last edit location) to\n// immediately follow the text we just replaced.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst endOfEdit = target.getRange(Word.RangeLocation.end);\nendOfEdit.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Update the SDK compatibility line:\n#   \"This sample is compatible with the Windows 10 Creators Update SDK (15063)\"\n# becomes:\n#   \"This sample is compatible with the Windows 10 Fall Creators Update SDK (16299)\"\n# (commit message: \"Updated for Windows 10 Fall Creators Update SDK (16299)\")\n#\n# Word records this as the user selecting just the version substring and\n# retyping it, which relocates the \"_GoBack\" last-edit bookmark to right\n# after the new text. We replicate both effects here.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: replace the SDK version text ---------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"Windows 10 Creators Update SDK (15063)\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Windows 10 Fall Creators Update SDK (16299)\", 2\n)\n\n# --- Step 2: drop the old \"_GoBack\" bookmark (Word's last-edit marker) ----\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- Step 3: re-create \"_GoBack\" immediately after the new text -----------\n# A collapsed range placed exactly at the end of a paragraph's text can't be\n# handed straight to Bookmarks.Add reliably, so stamp a temporary marker\n# right after the replacement, bookmark that span, then erase the marker -\n# the bookmark collapses in place, ending up right where the edit happened.\n$markerText = \"@@_GoBackMarker@@\"\n\n$insertFind = $d.Content.Find\n$insertFind.ClearFormatting()\n$insertFind.Replacement.ClearFormatting()\n$insertFind.Execute(\n    \"Windows 10 Fall Creators Update SDK (16299)\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Windows 10 Fall Creators Update SDK (16299)$markerText\", 2\n)\n\n$markerRange = $d.Content\n$markerFind = $markerRange.Find\n$markerFind.ClearFormatting()\n$markerFind.Text = $markerText\n$markerFind.Execute()\n\n$d.Bookmarks.Add(\"_GoBack\", $markerRange)\n$markerRange.Text = \"\"\n"}
